$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell G1: give it the same header style (s=2) as the other header cells ---
$ws.Range("G1").Value = "Note"
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats

# --- Header row values (Part Number / Mfr / Note shuffled across E/F/G) ---
$ws.Range("E1").Value = "Part Number"
$ws.Range("F1").Value = "Mfr"

# --- Row 2 (Item 1) ---
$ws.Range("D2").Value = "0.1uF 0402"
$ws.Range("F2").Value = "KEMET"

# --- Row 3 (Item 2) ---
$ws.Range("D3").Value = "1uF 0603"
$ws.Range("F3").Value = "KEMET"

# --- Row 4 (Item 3) ---
$ws.Range("D4").Value = "RGB LED 3210"
$ws.Range("F4").Value = "Adafruit"
$ws.Range("G4").Value = "mount upside-down"

# --- Row 5 (Item 4) ---
$ws.Range("D5").Value = "Conn 4 pos"
$ws.Range("F5").Value = "Samtec"

# --- Row 6 (Item 5): D6 already carries a quotePrefix-only style (s=1) with no value.
#     Preserve that style across the value write via a scratch-cell format round trip. ---
$scratch = $ws.Range("Z100")
$ws.Range("D6").Copy()
$scratch.PasteSpecial(-4122)
$ws.Range("D6").Value = "Conn 6 pos"
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4122)
$scratch.Clear()
$ws.Range("F6").Value = "Samtec"

# --- Row 7 (Item 6) ---
$ws.Range("D7").Value = "Card Edge Connector, 120"
$ws.Range("F7").Value = "Samtec"

# --- Row 8 (Item 7) ---
$ws.Range("D8").Value = "Header 30x2"
$ws.Range("F8").Value = "Samtec"

# --- Row 9 (Item 8) ---
$ws.Range("D9").Value = "PNP SOT323"
$ws.Range("F9").Value = "Nexperia"

# --- Row 10 (Item 9) ---
$ws.Range("F10").Value = "YAEGO"
$ws.Range("G10").Value = "DNP"

# --- Row 11 (Item 10) ---
$ws.Range("F11").Value = "YAEGO"

# --- Row 12 (Item 11) ---
$ws.Range("F12").Value = "YAEGO"

# --- Row 13 (Item 12) ---
$ws.Range("F13").Value = "YAEGO"

# --- Row 14 (Item 13) ---
$ws.Range("F14").Value = "YAEGO"

# --- Row 15 (Item 14) ---
$ws.Range("D15").Value = "SMT SW"
$ws.Range("F15").Value = "C&K"

# --- Column widths: new col F (Mfr) gets a fresh width, old col F width moves to new col G ---
$ws.Columns.Item(6).ColumnWidth = 25.7109375
$ws.Columns.Item(7).ColumnWidth = 40.28515625

# --- Selection restored on reopen ---
$ws.Range("F20").Select()
